$wb = $excel.ActiveWorkbook

# --- "Log" sheet: append new row 4 ---
$wsLog = $wb.Worksheets.Item("Log")
$wsLog.Range("A4").Value2 = 44646
$wsLog.Range("A4").NumberFormat = "d-mmm"
$wsLog.Range("B4").Value2 = "Replaced outliers with median values for the features which improved slightly the score" + [char]10 + "Checked the accuracy of the prediction using the real samples. Works fine"
$wsLog.Range("B4").WrapText = $true
$wsLog.Rows.Item(4).RowHeight = 28.8
$wsLog.Range("B7:B8").Select()

# --- "To Do" sheet: row 1 task is done (moved to Log), row 2's task is replaced ---
$wsToDo = $wb.Worksheets.Item("To Do")
$wsToDo.Rows.Item(1).ClearContents()
$wsToDo.Range("A2").Value2 = "Check what is wrong with estinmatedPrice it looks like the range is incorrect in streamlit slider"
$wsToDo.Columns.Item(1).ColumnWidth = 76.43
$wsToDo.Activate()
$wsToDo.Range("A8").Select()
